$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1575
$ws.Range("I38").Value = 150
$ws.Range("K38").Value = 450
$ws.Range("M38").Value = -78
$ws.Range("H61").Value = 1002.5
$ws.Range("I61").Value = 1002.5
$ws.Range("K61").Value = 3007.5
$ws.Range("M61").Value = -2835.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 583.6667
$ws.Range("I5").Value = 583.6667
$ws.Range("K5").Value = 583.6667
$ws.Range("M5").Value = -471.6667
$ws.Range("H32").Value = 1952389.6
$ws.Range("I32").Value = 1116.7693
$ws.Range("J32").Value = 6180147.5
$ws.Range("K32").Value = 1116.7693
$ws.Range("L32").Value = 6180147.5
$ws.Range("M32").Value = -829.7692999999999
$ws.Range("N32").Value = -6180721.5
$ws.Range("H45").Value = 1814.6
$ws.Range("I45").Value = 1705.2222
$ws.Range("K45").Value = 1705.2222
$ws.Range("M45").Value = -1328.2222
$ws.Range("H63").Value = 7014.2856
$ws.Range("J63").Value = 7014.2856
$ws.Range("L63").Value = 7014.2856
$ws.Range("N63").Value = -8386.285599999999
$ws.Range("H66").Value = 7014.2856
$ws.Range("J66").Value = 7014.2856
$ws.Range("L66").Value = 35071.428
$ws.Range("N66").Value = -41935.428
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H74").Value = 4700.7
$ws.Range("I74").Value = 4184.3335
$ws.Range("J74").Value = 6249.8
$ws.Range("K74").Value = 4184.3335
$ws.Range("L74").Value = 6249.8
$ws.Range("M74").Value = -3310.3335
$ws.Range("N74").Value = -7997.8
$ws.Range("H77").Value = 4700.7
$ws.Range("I77").Value = 4184.3335
$ws.Range("J77").Value = 6249.8
$ws.Range("K77").Value = 20921.6675
$ws.Range("L77").Value = 31249
$ws.Range("M77").Value = -16553.6675
$ws.Range("N77").Value = -39985
$ws.Range("H97").Value = 1804.9445
$ws.Range("I97").Value = 1483.6364
$ws.Range("J97").Value = 2309.8572
$ws.Range("K97").Value = 1483.6364
$ws.Range("L97").Value = 2309.8572
$ws.Range("M97").Value = -987.6364000000001
$ws.Range("N97").Value = -3301.8572
$ws.Range("H110").Value = 3220.913
$ws.Range("I110").Value = 1269.0714
$ws.Range("J110").Value = 6257.1113
$ws.Range("K110").Value = 1269.0714
$ws.Range("L110").Value = 6257.1113
$ws.Range("M110").Value = 775.9286
$ws.Range("N110").Value = -10347.1113
$ws.Range("H122").Value = 7046.5713
$ws.Range("I122").Value = 6604
$ws.Range("K122").Value = 19812
$ws.Range("M122").Value = -17362

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 583.6667
$ws.Range("I4").Value = 583.6667
$ws.Range("K4").Value = 583.6667
$ws.Range("M4").Value = -468.6667
$ws.Range("H20").Value = 4926836
$ws.Range("I20").Value = 6211863.5
$ws.Range("J20").Value = 895.8333
$ws.Range("K20").Value = 6211863.5
$ws.Range("L20").Value = 895.8333
$ws.Range("M20").Value = -6211616.5
$ws.Range("N20").Value = -1389.8333
$ws.Range("H94").Value = 5893.7715
$ws.Range("I94").Value = 2369.0386
$ws.Range("K94").Value = 2369.0386
$ws.Range("M94").Value = -1918.0386

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4303.6763
$ws.Range("I31").Value = 2042.375
$ws.Range("K31").Value = 2042.375
$ws.Range("M31").Value = -1747.375
$ws.Range("H34").Value = 4303.6763
$ws.Range("I34").Value = 2042.375
$ws.Range("K34").Value = 2042.375
$ws.Range("M34").Value = -1840.375
$ws.Range("H48").Value = 35045
$ws.Range("J48").Value = 35045
$ws.Range("L48").Value = 35045
$ws.Range("N48").Value = -35997
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H60").Value = 31832.834
$ws.Range("J60").Value = 36799.6
$ws.Range("L60").Value = 36799.6
$ws.Range("N60").Value = -37821.6
$ws.Range("H86").Value = 9561.866
$ws.Range("J86").Value = 13152.25
$ws.Range("L86").Value = 13152.25
$ws.Range("N86").Value = -15398.25
$ws.Range("H89").Value = 9561.866
$ws.Range("J89").Value = 13152.25
$ws.Range("L89").Value = 65761.25
$ws.Range("N89").Value = -76993.25
$ws.Range("H134").Value = 66676050
$ws.Range("I134").Value = 90916540
$ws.Range("J134").Value = 14678.25
$ws.Range("K134").Value = 272749620
$ws.Range("L134").Value = 44034.75
$ws.Range("M134").Value = -272747085
$ws.Range("N134").Value = -49104.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 86.15385000000001
$ws.Range("I2").Value = 64.5
$ws.Range("J2").Value = 90.09090999999999
$ws.Range("K2").Value = 387
$ws.Range("L2").Value = 540.5454599999999
$ws.Range("M2").Value = -274
$ws.Range("N2").Value = -766.5454599999999
$ws.Range("H6").Value = 3370
$ws.Range("J6").Value = 3370
$ws.Range("L6").Value = 10110
$ws.Range("N6").Value = -10336
$ws.Range("H38").Value = 69.45
$ws.Range("I38").Value = 62.636364
$ws.Range("J38").Value = 77.77778000000001
$ws.Range("K38").Value = 187.909092
$ws.Range("L38").Value = 233.33334
$ws.Range("M38").Value = 159.090908
$ws.Range("N38").Value = -927.33334
$ws.Range("H140").Value = 26424298
$ws.Range("I140").Value = 31863924
$ws.Range("K140").Value = 95591772
$ws.Range("M140").Value = -95586592
$ws.Range("H141").Value = 2737.8462
$ws.Range("I141").Value = 2737.8462
$ws.Range("K141").Value = 8213.5386
$ws.Range("M141").Value = -3033.5386

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 29795
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H50").Value = 29795
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()
$ws.Range("H70").Value = 14335
$ws.Range("I70").Value = 9587.333000000001
$ws.Range("J70").Value = 20032.2
$ws.Range("K70").Value = 9587.333000000001
$ws.Range("L70").Value = 20032.2
$ws.Range("M70").Value = -9317.333000000001
$ws.Range("N70").Value = -20572.2
$ws.Range("H73").Value = 14335
$ws.Range("I73").Value = 9587.333000000001
$ws.Range("J73").Value = 20032.2
$ws.Range("K73").Value = 9587.333000000001
$ws.Range("L73").Value = 20032.2
$ws.Range("M73").Value = -8651.333000000001
$ws.Range("N73").Value = -21904.2
$ws.Range("H80").Value = 58835320
$ws.Range("J80").Value = 15225.583
$ws.Range("L80").Value = 15225.583
$ws.Range("N80").Value = -17221.583
$ws.Range("H83").Value = 58835320
$ws.Range("J83").Value = 15225.583
$ws.Range("L83").Value = 76127.91500000001
$ws.Range("N83").Value = -86111.91500000001
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H132").Value = 45458744
$ws.Range("I132").Value = 90913060
$ws.Range("K132").Value = 272739180
$ws.Range("M132").Value = -272736650

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7297.5
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 7297.5
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 7297.5
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -7569.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H62").Value = 13109
$ws.Range("I62").Value = 10646.4
$ws.Range("K62").Value = 10646.4
$ws.Range("M62").Value = -10022.4
$ws.Range("H65").Value = 13109
$ws.Range("I65").Value = 10646.4
$ws.Range("K65").Value = 53232
$ws.Range("M65").Value = -50112
$ws.Range("H122").Value = 1802.3334
$ws.Range("I122").Value = 1866.6522
$ws.Range("K122").Value = 5599.9566
$ws.Range("M122").Value = -3149.9566
$ws.Range("H132").Value = 5541.1094
$ws.Range("I132").Value = 4677.981
$ws.Range("J132").Value = 9699.817999999999
$ws.Range("K132").Value = 14033.943
$ws.Range("L132").Value = 29099.454
$ws.Range("M132").Value = -11503.943
$ws.Range("N132").Value = -34159.454
$ws.Range("H136").Value = 11632521
$ws.Range("I136").Value = 18519324
$ws.Range("J136").Value = 11038.6875
$ws.Range("K136").Value = 55557972
$ws.Range("L136").Value = 33116.0625
$ws.Range("M136").Value = -55555422
$ws.Range("N136").Value = -38216.0625
